$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.826.17"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "1.812.62"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4317"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3707"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07243"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8658"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.01%  "

$ws.Range("D12").Value = "1.953.00"
$ws.Range("E12").Value = "  +4.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.639"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.31%  "

$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06921"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008919"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.96%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").Value = "26.863.41"
$ws.Range("E21").Value = "  -1.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.209"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "

$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").Value = "2.167.96"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.870"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.230"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.895"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7581"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.171"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.03%  "

$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("E35").Value = "  -1.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  +5.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05240"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01927"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5085"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1650"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.689"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.564"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.305"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.655"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4559"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06285"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.812"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.93%  "
